$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1084
$ws1.Range("F5").Value = 412
$ws1.Range("F7").Value = 541
$ws1.Range("F8").Value = 63
$ws1.Range("F9").Value = 6731
$ws1.Range("F10").Value = 151
$ws1.Range("F15").Value = 1082
$ws1.Range("F16").Value = 16090
$ws1.Range("F18").Value = 36
$ws1.Range("F22").Value = 11296
$ws1.Range("F23").Value = 6
$ws1.Range("F24").Value = 903
$ws1.Range("F28").Value = 41
$ws1.Range("F29").Value = 33
$ws1.Range("F30").Value = 317
$ws1.Range("F31").Value = 137
$ws1.Range("F32").Value = 5217

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1084
$ws4.Range("F5").Value = 412
$ws4.Range("F7").Value = 541
$ws4.Range("F9").Value = 63
$ws4.Range("F10").Value = 6731
$ws4.Range("F11").Value = 151
$ws4.Range("F15").Value = 0
$ws4.Range("F17").Value = 1082
$ws4.Range("F18").Value = 16090
$ws4.Range("F20").Value = 36
$ws4.Range("F26").Value = 11296
$ws4.Range("F27").Value = 6
$ws4.Range("F28").Value = 903
$ws4.Range("F32").Value = 41
$ws4.Range("F33").Value = 33
$ws4.Range("F34").Value = 317
$ws4.Range("F35").Value = 137
$ws4.Range("F36").Value = 5217
